## added home timezone to tournament info
#
# The "Tournament" sheet holds a key/value info table (the "tournament"
# ListObject). Insert a new "timezone" / "Europe/Moscow" row right after
# the existing "location" row (row 4) and before the "venue.1" row
# (previously row 5), shifting the venue rows down by one. Then make sure
# the table/ListObject grows to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")

# Insert a blank worksheet row above the old row 5 ("venue.1"), pushing
# every following row (venues, table formatting, etc.) down by one.
$ws.Rows.Item(5).Insert()

# Populate the new row with the timezone key/value pair.
$ws.Cells.Item(5, 1).Value = "timezone"
$ws.Cells.Item(5, 2).Value = "Europe/Moscow"

# Grow the "tournament" table so its range/autofilter include the new row.
$tbl = $ws.ListObjects.Item(1)
$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count
$lastCol = $tbl.Range.Column + $tbl.Range.Columns.Count - 1
$newRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl.Resize($newRange) | Out-Null

# Bring the Tournament sheet to the front and select C5, matching the
# author's final cursor position after adding the new row.
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null
